$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F ("dSF") updates per row
$updates = @{
    3  = -2
    4  = -2
    6  = 2
    7  = 10
    8  = 3
    10 = -1
    11 = -1
    12 = 1
    14 = 0
    15 = 1
    16 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
